$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.040.94"
$ws.Range("E2").Value = "  +2.24%  "
$ws.Range("D3").Value = "2.291.72"
$ws.Range("E3").Value = "  +3.20%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("E6").Value = "  +2.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.94"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.18%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.646"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0984"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "59.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("E13").Value = "  +5.04%  "
$ws.Range("E14").Value = "  +1.93%  "
$ws.Range("D15").Value = "2.636.02"
$ws.Range("E15").Value = "  +3.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.874"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "2.307.68"
$ws.Range("E18").Value = "  +3.14%  "
$ws.Range("D19").Value = "42.965.43"
$ws.Range("E19").Value = "  +2.25%  "
$ws.Range("E20").Value = "  +4.69%  "
$ws.Range("E21").Value = "  +2.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.51%  "
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.50%  "
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "166.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.128"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0821"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +16.12%  "
$ws.Range("E37").Value = "  +3.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.09%  "
$ws.Range("E39").Value = "  +3.28%  "
$ws.Range("E40").Value = "  -1.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +17.35%  "
$ws.Range("E42").Value = "  +4.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.217"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.03%  "
$ws.Range("E48").Value = "  +3.51%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("E50").Value = "  +1.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "99.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.42%  "
